$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.701699733734131
$ws.Range("B1").Value = 5.67929744720459
$ws.Range("C1").Value = 5.355366706848145
$ws.Range("D1").Value = 9.202488899230957
$ws.Range("E1").Value = 7.538563251495361
